$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kotik strut input")

$ws.Range("P2").Value = 686
$ws.Range("P3").Value = 686
$ws.Range("P4").Value = 1164
$ws.Range("P5").Value = 1164
$ws.Range("P6").Value = 1962
$ws.Range("P7").Value = 1962
$ws.Range("P8").Value = 1724
$ws.Range("P9").Value = 1724
$ws.Range("P10").Value = 960
$ws.Range("P11").Value = 960

$wb.Save()
